$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text for the two note cells
$text33 = "0.ZDK比较特殊：可以设置2~128大小字体，(#``O′未加粗"
$text36 = "如果想用其他大小，且加粗字体，就得生成对应字库，且数据变量的‘Y方向点阵数’要和新的字库封装大小一致"

# Write the values
$ws.Range("A33").Value = $text33
$ws.Range("A36").Value = $text36

# Format: font size 28, red color, font name 等线 (apply to each cell individually,
# so intervening empty rows/cells are not touched)
$font33 = $ws.Range("A33").Font
$font33.Name = "等线"
$font33.Size = 28
$font33.Color = 255

$font36 = $ws.Range("A36").Font
$font36.Name = "等线"
$font36.Size = 28
$font36.Color = 255

# Row heights
$ws.Rows.Item(33).RowHeight = 35.25
$ws.Rows.Item(36).RowHeight = 35.25

# Selection
$ws.Range("Z43").Select()
